# Auto-generated Excel COM-interop script implementing the target diff:
# - rename sheet, update selection, repopulate collector rows 2-18 across
#   columns A-L with updated figures from the refreshed report export
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force the value to be stored as text even when it looks numeric
    # (e.g. "2.37", "4,300,290.00"), matching the source report export.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
$ws.Range("A2").Value = "Riska Nurlita"
$ws.Range("B2").Value = "Hansyah_S2l"
$ws.Range("C2").Value = "S2"
$ws.Range("D2").Value = 6
Set-TextCell $ws.Range("E2") "4,300,290.00"
Set-TextCell $ws.Range("F2") "181,161,937.00"
Set-TextCell $ws.Range("G2") "2.37"
$ws.Range("H2").Value = 9
$ws.Range("I2").Value = 16
$ws.Range("J2").Value = 1
Set-TextCell $ws.Range("K2") "8.06"
Set-TextCell $ws.Range("L2") "6.25"

# Row 3
$ws.Range("A3").Value = "Debora Retima Sihombing"
$ws.Range("B3").Value = "Hansyah_S2l"
$ws.Range("C3").Value = "S2"
$ws.Range("D3").Value = 6
Set-TextCell $ws.Range("E3") "1,172,186.00"
Set-TextCell $ws.Range("F3") "163,604,043.00"
Set-TextCell $ws.Range("G3") "0.72"
$ws.Range("H3").Value = 320
$ws.Range("I3").Value = 17
$ws.Range("J3").Value = 0
Set-TextCell $ws.Range("K3") "0.00"
Set-TextCell $ws.Range("L3") "0.00"

# Row 4
$ws.Range("A4").Value = "Erlangga Hutama"
$ws.Range("B4").Value = "Hansyah_S2l"
$ws.Range("C4").Value = "S2"
$ws.Range("D4").Value = 11
Set-TextCell $ws.Range("E4") "3,362,575.00"
Set-TextCell $ws.Range("F4") "138,885,636.00"
Set-TextCell $ws.Range("G4") "2.42"
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 17
$ws.Range("J4").Value = 1
Set-TextCell $ws.Range("K4") "3.66"
Set-TextCell $ws.Range("L4") "5.88"

# Row 5
$ws.Range("A5").Value = "Romli"
$ws.Range("B5").Value = "Hansyah_S2l"
$ws.Range("C5").Value = "S2"
$ws.Range("D5").Value = 2
Set-TextCell $ws.Range("E5") "431,131.00"
Set-TextCell $ws.Range("F5") "138,503,407.00"
Set-TextCell $ws.Range("G5") "0.31"
$ws.Range("H5").Value = 217
$ws.Range("I5").Value = 17
$ws.Range("J5").Value = 0
Set-TextCell $ws.Range("K5") "1.41"
Set-TextCell $ws.Range("L5") "0.00"

# Row 6
$ws.Range("A6").Value = "Yandi Nugraha"
$ws.Range("B6").Value = "Hansyah_S2l"
$ws.Range("C6").Value = "S2"
$ws.Range("D6").Value = 5
Set-TextCell $ws.Range("E6") "2,653,033.00"
Set-TextCell $ws.Range("F6") "132,519,967.00"
Set-TextCell $ws.Range("G6") "2.00"
$ws.Range("H6").Value = 16
$ws.Range("I6").Value = 17
$ws.Range("J6").Value = 0
Set-TextCell $ws.Range("K6") "0.00"
Set-TextCell $ws.Range("L6") "0.00"

# Row 7
$ws.Range("A7").Value = "Sucika Wardani"
$ws.Range("B7").Value = "Hansyah_S2l"
$ws.Range("C7").Value = "S2"
$ws.Range("D7").Value = 5
Set-TextCell $ws.Range("E7") "488,364.00"
Set-TextCell $ws.Range("F7") "175,090,870.00"
Set-TextCell $ws.Range("G7") "0.28"
$ws.Range("H7").Value = 18
$ws.Range("I7").Value = 17
$ws.Range("J7").Value = 0
Set-TextCell $ws.Range("K7") "0.00"
Set-TextCell $ws.Range("L7") "0.00"

# Row 8
$ws.Range("A8").Value = "Nuraini"
$ws.Range("B8").Value = "Hansyah_S2l"
$ws.Range("C8").Value = "S2"
$ws.Range("D8").Value = 2
Set-TextCell $ws.Range("E8") "3,393,304.00"
Set-TextCell $ws.Range("F8") "124,985,498.00"
Set-TextCell $ws.Range("G8") "2.71"
$ws.Range("H8").Value = 34
$ws.Range("I8").Value = 16
$ws.Range("J8").Value = 0
Set-TextCell $ws.Range("K8") "0.00"
Set-TextCell $ws.Range("L8") "0.00"

# Row 9
$ws.Range("A9").Value = "Annisa Putri Restu"
$ws.Range("B9").Value = "Hansyah_S2l"
$ws.Range("C9").Value = "S2"
$ws.Range("D9").Value = 6
Set-TextCell $ws.Range("E9") "4,462,496.00"
Set-TextCell $ws.Range("F9") "182,088,660.00"
Set-TextCell $ws.Range("G9") "2.45"
$ws.Range("H9").Value = 260
$ws.Range("I9").Value = 16
$ws.Range("J9").Value = 0
Set-TextCell $ws.Range("K9") "0.00"
Set-TextCell $ws.Range("L9") "0.00"

# Row 10
$ws.Range("A10").Value = "Azizah Rahmawati"
$ws.Range("B10").Value = "Hansyah_S2l"
$ws.Range("C10").Value = "S2"
$ws.Range("D10").Value = 10
Set-TextCell $ws.Range("E10") "3,395,577.00"
Set-TextCell $ws.Range("F10") "171,821,552.00"
Set-TextCell $ws.Range("G10") "1.98"
$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 16
$ws.Range("J10").Value = 1
Set-TextCell $ws.Range("K10") "0.80"
Set-TextCell $ws.Range("L10") "6.25"

# Row 11
$ws.Range("A11").Value = "Aldi Taufik"
$ws.Range("B11").Value = "Hansyah_S2l"
$ws.Range("C11").Value = "S2"
$ws.Range("D11").Value = 1
Set-TextCell $ws.Range("E11") "1,453,709.00"
Set-TextCell $ws.Range("F11") "153,773,268.00"
Set-TextCell $ws.Range("G11") "0.95"
$ws.Range("H11").Value = 13
$ws.Range("I11").Value = 16
$ws.Range("J11").Value = 1
Set-TextCell $ws.Range("K11") "13.23"
Set-TextCell $ws.Range("L11") "6.25"

# Row 12
$ws.Range("A12").Value = "Axl Wicaksono"
$ws.Range("B12").Value = "Hansyah_S2l"
$ws.Range("C12").Value = "S2"
$ws.Range("D12").Value = 1
Set-TextCell $ws.Range("E12") "90,243.00"
Set-TextCell $ws.Range("F12") "149,436,886.00"
Set-TextCell $ws.Range("G12") "0.06"
$ws.Range("H12").Value = 154
$ws.Range("I12").Value = 16
$ws.Range("J12").Value = 0
Set-TextCell $ws.Range("K12") "0.00"
Set-TextCell $ws.Range("L12") "0.00"

# Row 13
$ws.Range("A13").Value = "Ridhoi Berkat Zebua"
$ws.Range("B13").Value = "Hansyah_S2l"
$ws.Range("C13").Value = "S2"
$ws.Range("D13").Value = 5
Set-TextCell $ws.Range("E13") "1,134,065.00"
Set-TextCell $ws.Range("F13") "165,151,431.00"
Set-TextCell $ws.Range("G13") "0.69"
$ws.Range("H13").Value = 454
$ws.Range("I13").Value = 16
$ws.Range("J13").Value = 1
Set-TextCell $ws.Range("K13") "4.03"
Set-TextCell $ws.Range("L13") "6.25"

# Row 14
$ws.Range("A14").Value = "Fadilah Damayanti"
$ws.Range("B14").Value = "Hansyah_S2l"
$ws.Range("C14").Value = "S2"
$ws.Range("D14").Value = 1
Set-TextCell $ws.Range("E14") "50,000.00"
Set-TextCell $ws.Range("F14") "136,709,694.00"
Set-TextCell $ws.Range("G14") "0.04"
$ws.Range("H14").Value = 83
$ws.Range("I14").Value = 16
$ws.Range("J14").Value = 0
Set-TextCell $ws.Range("K14") "0.00"
Set-TextCell $ws.Range("L14") "0.00"

# Row 15
$ws.Range("A15").Value = "Adistira Winditya P"
$ws.Range("B15").Value = "Hansyah_S2l"
$ws.Range("C15").Value = "S2"
$ws.Range("D15").Value = 1
Set-TextCell $ws.Range("E15") "379,713.00"
Set-TextCell $ws.Range("F15") "164,878,289.00"
Set-TextCell $ws.Range("G15") "0.23"
$ws.Range("H15").Value = 5
$ws.Range("I15").Value = 17
$ws.Range("J15").Value = 0
Set-TextCell $ws.Range("K15") "0.00"
Set-TextCell $ws.Range("L15") "0.00"

# Row 16
$ws.Range("A16").Value = "Erick Ervan Dewanggga"
$ws.Range("B16").Value = "Hansyah_S2l"
$ws.Range("C16").Value = "S2"
$ws.Range("D16").Value = 1
Set-TextCell $ws.Range("E16") "200,000.00"
Set-TextCell $ws.Range("F16") "168,753,581.00"
Set-TextCell $ws.Range("G16") "0.12"
$ws.Range("H16").Value = 27
$ws.Range("I16").Value = 16
$ws.Range("J16").Value = 0
Set-TextCell $ws.Range("K16") "0.00"
Set-TextCell $ws.Range("L16") "0.00"

# Row 17
$ws.Range("A17").Value = "Nur Halim"
$ws.Range("B17").Value = "Hansyah_S2l"
$ws.Range("C17").Value = "S2"
$ws.Range("D17").Value = 2
Set-TextCell $ws.Range("E17") "1,027,238.00"
Set-TextCell $ws.Range("F17") "129,576,507.00"
Set-TextCell $ws.Range("G17") "0.79"
$ws.Range("H17").Value = 74
$ws.Range("I17").Value = 16
$ws.Range("J17").Value = 0
Set-TextCell $ws.Range("K17") "0.00"
Set-TextCell $ws.Range("L17") "0.00"

# Row 18
$ws.Range("A18").Value = "Wasti Feronika Sihombing"
$ws.Range("B18").Value = "Hansyah_S2l"
$ws.Range("C18").Value = "S2"
$ws.Range("D18").Value = 1
Set-TextCell $ws.Range("E18") "557,294.00"
Set-TextCell $ws.Range("F18") "139,820,937.00"
Set-TextCell $ws.Range("G18") "0.40"
$ws.Range("H18").Value = 116
$ws.Range("I18").Value = 17
$ws.Range("J18").Value = 0
Set-TextCell $ws.Range("K18") "0.00"
Set-TextCell $ws.Range("L18") "0.00"

# Rename the worksheet tab (duplicated upload -> "(3)" suffix)
$ws.Name = "repayment_20250925_20250925 (3)"

# Update the active selection to match the refreshed data range
$ws.Range("A2:A18").Select()
